$d = $word.ActiveDocument
$section = $d.Sections.First
$footer = $section.Footers.Item(1)
$footer.Range.Tables.Add($footer.Range, 1, 3)
